$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.461.13'
$ws.Range('E2').Value = '  +1.00%  '
$ws.Range('D3').Value = '1.878.00'
$ws.Range('E3').Value = '  +0.85%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '246.69'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +5.46%  '
$ws.Range('E6').Value = '  -0.12%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4763'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +1.82%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2902'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +1.78%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06525'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +0.81%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.90'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +4.25%  '
$ws.Range('E11').Value = '  -0.28%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '97.30'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +4.07%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.7383'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +8.41%  '
$ws.Range('D14').Value = '1.878.45'
$ws.Range('E14').Value = '  +0.71%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.131'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +1.79%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '272.86'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +2.00%  '
$ws.Range('D17').Value = '30.446.02'
$ws.Range('E17').Value = '  +0.99%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.62'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +2.32%  '
$ws.Range('E19').Value = '  -0.13%  '
$ws.Range('D21').Value = '2.121.17'
$ws.Range('E21').Value = '  +0.29%  '
$ws.Range('E22').Value = '  -0.05%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.254'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +2.33%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.186'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +1.43%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.336'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.04%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '164.08'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -0.53%  '
$ws.Range('E27').Value = '  +2.01%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.940'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +2.89%  '
$ws.Range('E29').Value = '  +0.57%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.09948'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +0.04%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.518'
$ws.Range('D31').ClearFormats()
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.313'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +2.27%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.069'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +1.94%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.04785'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +2.54%  '
$ws.Range('E35').Value = '  +0.99%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7006'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +1.63%  '
$ws.Range('E37').Value = '  -0.01%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01866'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +2.03%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.730'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -1.12%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.341'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +0.62%  '
$ws.Range('B41').Value = 'Aave'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '70.91'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -0.38%  '
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.948'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +3.48%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.4203'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +3.91%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9997'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -0.10%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8372'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +0.79%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '102.83'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +0.78%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.255'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +1.26%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.092'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +2.15%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '35.63'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +4.78%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '927.98'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -0.32%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05646'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +1.30%  '
